$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '26.413.68'
$ws.Range('D3').Value = '1.618.79'
$ws.Range('E3').Value = '  +1.52%  '
$ws.Range('E4').Value = '  -0.01%  '
$c = $ws.Range('D5')
$c.NumberFormat = "@"
$c.Value = '212.99'
$c.ClearFormats()
$ws.Range('E5').Value = '  +0.08%  '
$ws.Range('E6').Value = '  +0.37%  '
$ws.Range('E7').Value = '  +0.02%  '
$ws.Range('E8').Value = '  +0.22%  '
$ws.Range('E9').Value = '  +0.32%  '
$c = $ws.Range('D10')
$c.NumberFormat = "@"
$c.Value = '19.19'
$c.ClearFormats()
$c = $ws.Range('D11')
$c.NumberFormat = "@"
$c.Value = '0.0849'
$c.ClearFormats()
$ws.Range('E11').Value = '  -0.39%  '
$ws.Range('D12').Value = '1.845.59'
$ws.Range('E12').Value = '  +1.55%  '
$ws.Range('D13').Value = '1.614.55'
$ws.Range('E13').Value = '  +0.76%  '
$ws.Range('E14').Value = '  +0.23%  '
$ws.Range('E15').Value = '  +0.30%  '
$c = $ws.Range('D16')
$c.NumberFormat = "@"
$c.Value = '63.96'
$c.ClearFormats()
$c = $ws.Range('D17')
$c.NumberFormat = "@"
$c.Value = '237.53'
$c.ClearFormats()
$ws.Range('E17').Value = '  +10.19%  '
$ws.Range('D18').Value = '26.423.48'
$ws.Range('E18').Value = '  +0.59%  '
$c = $ws.Range('D19')
$c.NumberFormat = "@"
$c.Value = '7.83'
$c.ClearFormats()
$ws.Range('E19').Value = '  +5.10%  '
$ws.Range('E20').Value = '  +0.43%  '
$ws.Range('E21').Value = '  +0.01%  '
$ws.Range('E22').Value = '  -0.01%  '
$c = $ws.Range('D23')
$c.NumberFormat = "@"
$c.Value = '9.11'
$c.ClearFormats()
$ws.Range('E23').Value = '  +1.28%  '
$ws.Range('E24').Value = '  +2.80%  '
$c = $ws.Range('D25')
$c.NumberFormat = "@"
$c.Value = '147.07'
$c.ClearFormats()
$ws.Range('E25').Value = '  +1.57%  '
$ws.Range('E26').Value = '  -0.07%  '
$c = $ws.Range('D27')
$c.NumberFormat = "@"
$c.Value = '7.05'
$c.ClearFormats()
$ws.Range('E27').Value = '  +1.22%  '
$ws.Range('E28').Value = '  +0.62%  '
$c = $ws.Range('D29')
$c.NumberFormat = "@"
$c.Value = '15.53'
$c.ClearFormats()
$ws.Range('E29').Value = '  +2.61%  '
$ws.Range('E30').Value = '  +0.32%  '
$ws.Range('E31').Value = '  -0.05%  '
$ws.Range('D32').Value = '1.526.52'
$ws.Range('E32').Value = '  +6.77%  '
$c = $ws.Range('D33')
$c.NumberFormat = "@"
$c.Value = '3.25'
$c.ClearFormats()
$ws.Range('E33').Value = '  +1.50%  '
$ws.Range('E34').Value = '  +0.35%  '
$ws.Range('E35').Value = '  +4.18%  '
$ws.Range('E36').Value = '  +0.37%  '
$c = $ws.Range('D37')
$c.NumberFormat = "@"
$c.Value = '0.570'
$c.ClearFormats()
$ws.Range('E37').Value = '  +1.64%  '
$ws.Range('E38').Value = '  +0.52%  '
$c = $ws.Range('D39')
$c.NumberFormat = "@"
$c.Value = '0.833'
$c.ClearFormats()
$ws.Range('E39').Value = '  +0.90%  '
$c = $ws.Range('D40')
$c.NumberFormat = "@"
$c.Value = '5.93'
$c.ClearFormats()
$ws.Range('E40').Value = '  +2.74%  '
$ws.Range('E41').Value = '  -0.05%  '
$ws.Range('E42').Value = '  +1.73%  '
$ws.Range('D43').Value = '1.757.44'
$ws.Range('E43').Value = '  +1.53%  '
$c = $ws.Range('D44')
$c.NumberFormat = "@"
$c.Value = '0.764'
$c.ClearFormats()
$ws.Range('E44').Value = '  +0.92%  '
$c = $ws.Range('D45')
$c.NumberFormat = "@"
$c.Value = '61.73'
$c.ClearFormats()
$ws.Range('E45').Value = '  +1.33%  '
$ws.Range('E46').Value = '  -0.17%  '
$c = $ws.Range('D47')
$c.NumberFormat = "@"
$c.Value = '90.83'
$c.ClearFormats()
$ws.Range('E47').Value = '  +4.73%  '
$ws.Range('E48').Value = '  +2.00%  '
$ws.Range('E49').Value = '  +0.30%  '
$c = $ws.Range('D50')
$c.NumberFormat = "@"
$c.Value = '0.0964'
$c.ClearFormats()
$ws.Range('E50').Value = '  +1.24%  '
$c = $ws.Range('D51')
$c.NumberFormat = "@"
$c.Value = '7.49'
$c.ClearFormats()
$ws.Range('E51').Value = '  +0.79%  '
